$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 20534
$ws.Range("E2").Value = 699
$ws.Range("F2").Value = 699
$ws.Range("G2").Value = 803
$ws.Range("H2").Value = 638
$ws.Range("I2").Value = 641
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 14404
$ws.Range("L2").Value = 5983
$ws.Range("M2").Value = 8421
$ws.Range("N2").Value = 8361
$ws.Range("O2").Value = 60
$ws.Range("P2").Value = 410
$ws.Range("Q2").Value = 775
$ws.Range("R2").Value = -1725
$ws.Range("S2").Value = 794
$ws.Range("T2").Value = 1638
$ws.Range("U2").Value = -863
$ws.Range("V2").Value = 3856
$ws.Range("W2").Value = 3.4
$ws.Range("X2").Value = 3.11
$ws.Range("Y2").Value = 7.92
$ws.Range("Z2").Value = 4.57
$ws.Range("AA2").Value = 71.05
$ws.Range("AB2").Value = 1948.52
$ws.Range("AC2").Value = 9864
$ws.Range("AD2").Value = 7.19
$ws.Range("AE2").Value = 135384
$ws.Range("AF2").Value = 0.52
$ws.Range("AG2").Value = 1500
$ws.Range("AH2").Value = 2.12
$ws.Range("AI2").Value = 14.45
$ws.Range("AJ2").Value = 6500000

# Row 3
$ws.Range("D3").Value = 17270
$ws.Range("E3").Value = 2712
$ws.Range("F3").Value = 2712
$ws.Range("G3").Value = 2687
$ws.Range("H3").Value = 2008
$ws.Range("I3").Value = 2005
$ws.Range("J3").Value = 3
$ws.Range("K3").Value = 15291
$ws.Range("L3").Value = 4935
$ws.Range("M3").Value = 10356
$ws.Range("N3").Value = 10293
$ws.Range("O3").Value = 62
$ws.Range("P3").Value = 410
$ws.Range("Q3").Value = 3237
$ws.Range("R3").Value = -1138
$ws.Range("S3").Value = -1359
$ws.Range("T3").Value = 844
$ws.Range("U3").Value = 2394
$ws.Range("V3").Value = 2603
$ws.Range("W3").Value = 15.7
$ws.Range("X3").Value = 11.63
$ws.Range("Y3").Value = 21.5
$ws.Range("Z3").Value = 13.52
$ws.Range("AA3").Value = 47.65
$ws.Range("AB3").Value = 2419.32
$ws.Range("AC3").Value = 30846
$ws.Range("AD3").Value = 5.58
$ws.Range("AE3").Value = 166666
$ws.Range("AF3").Value = 1.03
$ws.Range("AG3").Value = 3000
$ws.Range("AH3").Value = 1.74
$ws.Range("AI3").Value = 9.24
$ws.Range("AJ3").Value = 6500000

# Row 4
$ws.Range("D4").Value = 15964
$ws.Range("E4").Value = 3430
$ws.Range("F4").Value = 3430
$ws.Range("G4").Value = 3602
$ws.Range("H4").Value = 2729
$ws.Range("I4").Value = 2721
$ws.Range("J4").Value = 8
$ws.Range("K4").Value = 17636
$ws.Range("L4").Value = 4706
$ws.Range("M4").Value = 12931
$ws.Range("N4").Value = 12861
$ws.Range("O4").Value = 70
$ws.Range("P4").Value = 410
$ws.Range("Q4").Value = 3226
$ws.Range("R4").Value = -2780
$ws.Range("S4").Value = -919
$ws.Range("T4").Value = 2627
$ws.Range("U4").Value = 598
$ws.Range("V4").Value = 1865
$ws.Range("W4").Value = 21.49
$ws.Range("X4").Value = 17.09
$ws.Range("Y4").Value = 23.51
$ws.Range("Z4").Value = 16.58
$ws.Range("AA4").Value = 36.39
$ws.Range("AB4").Value = 3045.72
$ws.Range("AC4").Value = 41867
$ws.Range("AD4").Value = 6.58
$ws.Range("AE4").Value = 208231
$ws.Range("AF4").Value = 1.32
$ws.Range("AG4").Value = 4000
$ws.Range("AH4").Value = 1.45
$ws.Range("AI4").Value = 9.08
$ws.Range("AJ4").Value = 6500000

# Row 5
$ws.Range("D5").Value = 17794
$ws.Range("E5").Value = 2842
$ws.Range("F5").Value = 2842
$ws.Range("G5").Value = 2801
$ws.Range("H5").Value = 2145
$ws.Range("I5").Value = 2143
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 19706
$ws.Range("L5").Value = 4974
$ws.Range("M5").Value = 14732
$ws.Range("N5").Value = 14732
$ws.Range("P5").Value = 410
$ws.Range("Q5").Value = 2279
$ws.Range("R5").Value = -3064
$ws.Range("S5").Value = 11
$ws.Range("T5").Value = 3183
$ws.Range("U5").Value = -905
$ws.Range("V5").Value = 2215
$ws.Range("W5").Value = 15.97
$ws.Range("X5").Value = 12.05
$ws.Range("Y5").Value = 15.53
$ws.Range("Z5").Value = 11.49
$ws.Range("AA5").Value = 33.76
$ws.Range("AB5").Value = 3507.37
$ws.Range("AC5").Value = 32972
$ws.Range("AD5").Value = 7.99
$ws.Range("AE5").Value = 238533
$ws.Range("AF5").Value = 1.1
$ws.Range("AG5").Value = 4000
$ws.Range("AH5").Value = 1.52
$ws.Range("AI5").Value = 11.53
$ws.Range("AJ5").Value = 6500000
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 25540
$ws.Range("E6").Value = 3188
$ws.Range("F6").Value = 3188
$ws.Range("G6").Value = 3364
$ws.Range("H6").Value = 2574
$ws.Range("I6").Value = 2574
$ws.Range("K6").Value = 19839
$ws.Range("L6").Value = 2828
$ws.Range("M6").Value = 17011
$ws.Range("N6").Value = 17011
$ws.Range("P6").Value = 410
$ws.Range("Q6").Value = 3509
$ws.Range("R6").Value = -1166
$ws.Range("S6").Value = -1902
$ws.Range("T6").Value = 1258
$ws.Range("U6").Value = 2251
$ws.Range("V6").Value = 560
$ws.Range("W6").Value = 12.48
$ws.Range("X6").Value = 10.08
$ws.Range("Y6").Value = 16.22
$ws.Range("Z6").Value = 13.02
$ws.Range("AA6").Value = 16.63
$ws.Range("AB6").Value = 4063.24
$ws.Range("AC6").Value = 39602
$ws.Range("AD6").Value = 3.88
$ws.Range("AE6").Value = 275432
$ws.Range("AF6").Value = 0.5600000000000001
$ws.Range("AG6").Value = 4000
$ws.Range("AH6").Value = 2.61
$ws.Range("AI6").Value = 9.6
$ws.Range("AJ6").Value = 6500000

# Row 7
$ws.Range("D7").Value = 20451
$ws.Range("E7").Value = 1295
$ws.Range("G7").Value = 1458
$ws.Range("H7").Value = 1173
$ws.Range("I7").Value = 1173
$ws.Range("K7").Value = 21108
$ws.Range("L7").Value = 3217
$ws.Range("M7").Value = 17891
$ws.Range("N7").Value = 17895
$ws.Range("P7").Value = 410
$ws.Range("Q7").Value = 2438
$ws.Range("R7").Value = -1504
$ws.Range("S7").Value = -150
$ws.Range("T7").Value = 1575
$ws.Range("U7").Value = 980
$ws.Range("W7").Value = 6.33
$ws.Range("X7").Value = 5.73
$ws.Range("Y7").Value = 6.72
$ws.Range("Z7").Value = 5.73
$ws.Range("AA7").Value = 17.98
$ws.Range("AC7").Value = 18042
$ws.Range("AD7").Value = 5.99
$ws.Range("AE7").Value = 289746
$ws.Range("AF7").Value = 0.37
$ws.Range("AG7").Value = 3600
$ws.Range("AH7").Value = 3.33
$ws.Range("AI7").Value = 19.95

# Row 8
$ws.Range("D8").Value = 22674
$ws.Range("E8").Value = 1380
$ws.Range("G8").Value = 1547
$ws.Range("H8").Value = 1256
$ws.Range("I8").Value = 1191
$ws.Range("K8").Value = 22250
$ws.Range("L8").Value = 3374
$ws.Range("M8").Value = 18876
$ws.Range("N8").Value = 18939
$ws.Range("P8").Value = 410
$ws.Range("Q8").Value = 2866
$ws.Range("R8").Value = -1351
$ws.Range("S8").Value = -159
$ws.Range("T8").Value = 1362
$ws.Range("U8").Value = 1137
$ws.Range("W8").Value = 6.09
$ws.Range("X8").Value = 5.54
$ws.Range("Y8").Value = 6.47
$ws.Range("Z8").Value = 5.79
$ws.Range("AA8").Value = 17.88
$ws.Range("AC8").Value = 18327
$ws.Range("AD8").Value = 5.1
$ws.Range("AE8").Value = 306650
$ws.Range("AF8").Value = 0.3
$ws.Range("AG8").Value = 3667
$ws.Range("AH8").Value = 3.92
$ws.Range("AI8").Value = 20.01

# Row 9
$ws.Range("D9").Value = 23102
$ws.Range("E9").Value = 1619
$ws.Range("G9").Value = 1811
$ws.Range("H9").Value = 1473
$ws.Range("I9").Value = 1398
$ws.Range("K9").Value = 23638
$ws.Range("L9").Value = 3526
$ws.Range("M9").Value = 20111
$ws.Range("N9").Value = 20256
$ws.Range("P9").Value = 410
$ws.Range("Q9").Value = 3209
$ws.Range("R9").Value = -1482
$ws.Range("S9").Value = -118
$ws.Range("T9").Value = 1285
$ws.Range("U9").Value = 1144
$ws.Range("W9").Value = 7.01
$ws.Range("X9").Value = 6.38
$ws.Range("Y9").Value = 7.13
$ws.Range("Z9").Value = 6.42
$ws.Range("AA9").Value = 17.53
$ws.Range("AC9").Value = 21510
$ws.Range("AD9").Value = 4.35
$ws.Range("AE9").Value = 327972
$ws.Range("AF9").Value = 0.29
$ws.Range("AG9").Value = 3708
$ws.Range("AH9").Value = 3.97
$ws.Range("AI9").Value = 17.24
